$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for existing rows (AgTests / AgPosit revisions)
$ws.Range("F310").Value = 79480
$ws.Range("F322").Value = 110376
$ws.Range("F324").Value = 248840
$ws.Range("F327").Value = 225531
$ws.Range("F330").Value = 72934
$ws.Range("G330").Value = 2090
$ws.Range("F337").Value = 105628
$ws.Range("F352").Value = 307957
$ws.Range("F359").Value = 321003
$ws.Range("F363").Value = 189558
$ws.Range("F364").Value = 168770
$ws.Range("F366").Value = 339256
$ws.Range("F370").Value = 180679
$ws.Range("F371").Value = 160321
$ws.Range("F372").Value = 178701
$ws.Range("F373").Value = 350754
$ws.Range("F375").Value = 350347
$ws.Range("F376").Value = 222844
$ws.Range("F378").Value = 157511
$ws.Range("F379").Value = 180843
$ws.Range("F383").Value = 222698
$ws.Range("F386").Value = 183269
$ws.Range("F390").Value = 220022
$ws.Range("F392").Value = 222030
$ws.Range("F398").Value = 300519
$ws.Range("F403").Value = 353973
$ws.Range("F405").Value = 174938
$ws.Range("F408").Value = 305867
$ws.Range("F415").Value = 308819
$ws.Range("F421").Value = 153358
$ws.Range("F422").Value = 298727
$ws.Range("F426").Value = 107418
$ws.Range("F427").Value = 90525
$ws.Range("F428").Value = 102547
$ws.Range("F429").Value = 178478
$ws.Range("F432").Value = 122697
$ws.Range("F436").Value = 145444
$ws.Range("F439").Value = 89309
$ws.Range("F442").Value = 70543
$ws.Range("F443").Value = 106932
$ws.Range("F447").Value = 67042
$ws.Range("F449").Value = 60012
$ws.Range("F453").Value = 70162
$ws.Range("F456").Value = 50375
$ws.Range("F463").Value = 46683
$ws.Range("F464").Value = 73614
$ws.Range("F467").Value = 52245
$ws.Range("F470").Value = 43531
$ws.Range("F471").Value = 66831
$ws.Range("F472").Value = 51790
$ws.Range("G472").Value = 23
$ws.Range("F473").Value = 39845
$ws.Range("F475").Value = 36643
$ws.Range("F476").Value = 37313
$ws.Range("F477").Value = 36962
$ws.Range("F478").Value = 54886
$ws.Range("F479").Value = 42545
$ws.Range("F480").Value = 33312
$ws.Range("F481").Value = 41340
$ws.Range("F482").Value = 36110
$ws.Range("F483").Value = 65174
$ws.Range("G483").Value = 38
$ws.Range("F484").Value = 8270
$ws.Range("F486").Value = 8871
$ws.Range("F487").Value = 6805
$ws.Range("F488").Value = 6047
$ws.Range("F489").Value = 12625
$ws.Range("F490").Value = 10682
$ws.Range("F492").Value = 13835
$ws.Range("F493").Value = 8117
$ws.Range("F494").Value = 6274
$ws.Range("F495").Value = 10181
$ws.Range("F496").Value = 8027
$ws.Range("F497").Value = 7521
$ws.Range("F498").Value = 8970
$ws.Range("F499").Value = 10627
$ws.Range("F500").Value = 7057
$ws.Range("F501").Value = 5607
$ws.Range("F502").Value = 9455
$ws.Range("F503").Value = 7052
$ws.Range("F504").Value = 6859
$ws.Range("G504").Value = 16

# New daily rows through 2021-07-25
$ws.Range("A505").Value = 44399
$ws.Range("B505").Value = 392259
$ws.Range("C505").Value = 7466
$ws.Range("D505").Value = 40
$ws.Range("E505").Value = 12534
$ws.Range("F505").Value = 7612
$ws.Range("G505").Value = 27
$ws.Range("A506").Value = 44400
$ws.Range("B506").Value = 392302
$ws.Range("C506").Value = 10540
$ws.Range("D506").Value = 43
$ws.Range("E506").Value = 12534
$ws.Range("F506").Value = 9138
$ws.Range("G506").Value = 5
$ws.Range("A507").Value = 44401
$ws.Range("B507").Value = 392348
$ws.Range("C507").Value = 5834
$ws.Range("D507").Value = 46
$ws.Range("E507").Value = 12534
$ws.Range("F507").Value = 5872
$ws.Range("G507").Value = 8
$ws.Range("A508").Value = 44402
$ws.Range("B508").Value = 392355
$ws.Range("C508").Value = 1432
$ws.Range("D508").Value = 7
$ws.Range("E508").Value = 12534
$ws.Range("F508").Value = 4079
$ws.Range("G508").Value = 11
